$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "314.77"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.62%"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "21"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "49.02"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "10.61%"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "21"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.290"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "3.64%"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "21"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07897"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-1.06%"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "21"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.591"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "2.80%"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "21"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.323"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "22.86%"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "21"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.617"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.19%"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "21"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1237"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-4.10%"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "21"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1951"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2.91%"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "21"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09453"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.81%"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "21"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04544"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "7.64%"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "21"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1048"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "1.21%"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "21"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001309"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.03%"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "21"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.04210"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.64%"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "21"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005836"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "2.31%"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "21"
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "HotbitToken"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.004187"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-8.63%"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "21"
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.343"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.95%"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "21"
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.467"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2.73%"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "21"
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3465"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "3.16%"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "21"
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "MCDex"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.055"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.71%"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "21"
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = "ProBitToken"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.1405"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "2.21%"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "21"
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = "ZBToken"
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.3074"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-1.59%"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "21"
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = "BitKan"
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001294"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.83%"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "21"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "1.54%"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "21"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-95.19%"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "21"
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "21"
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "21"
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "21"
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "21"
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "21"
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "21"
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "21"
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "21"
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "21"
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "21"
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "21"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02641"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-0.58%"
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "21"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05837"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "7.65%"
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "21"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01077"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "92.15%"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "21"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.008014"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "3.85%"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "21"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1444"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "2.46%"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "21"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007809"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "7.17%"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "21"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008656"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "3.07%"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "21"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3145"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "0.72%"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "21"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00007032"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "4.74%"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "21"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000754"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "1.69%"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "21"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-7.83%"
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "21"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.004024"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "1.69%"
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "21"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002112"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "1.69%"
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "21"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002012"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "1.69%"
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "21"
